$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the shared text "HALO" with distinct values reflecting detection/
# extraction/comparison of excel file states.
$ws.Range("F3").Value = "H"
$ws.Range("F2").Value = 2323321
$ws.Range("F4").Value = "H2"
$ws.Range("F5").Value = "H"
$ws.Range("F6").Value = "H2323"
$ws.Range("F7").Value = "HALO111"

$ws.Range("F10").Select()
